$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '43.616.60'
Set-TextValue 'E2' '  +1.12%  '
Set-TextValue 'D3' '2.413.80'
Set-TextValue 'E3' '  +2.61%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '306.73'
Set-TextValue 'E5' '  +1.35%  '
Set-TextValue 'D6' '97.16'
Set-TextValue 'E6' '  +1.69%  '
Set-TextValue 'E7' '  +0.63%  '
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'D9' '0.492'
Set-TextValue 'E9' '  -1.14%  '
Set-TextValue 'D10' '35.12'
Set-TextValue 'E10' '  +3.13%  '
Set-TextValue 'D11' '0.0799'
Set-TextValue 'E11' '  +1.36%  '
Set-TextValue 'E12' '  +2.43%  '
Set-TextValue 'D13' '18.47'
Set-TextValue 'E13' '  -1.26%  '
Set-TextValue 'D14' '6.89'
Set-TextValue 'E14' '  +2.55%  '
Set-TextValue 'D15' '2.781.46'
Set-TextValue 'E15' '  +2.09%  '
Set-TextValue 'D16' '2.395.87'
Set-TextValue 'E16' '  +0.67%  '
Set-TextValue 'D17' '0.826'
Set-TextValue 'E17' '  +3.93%  '
Set-TextValue 'D18' '43.603.71'
Set-TextValue 'E18' '  +1.09%  '
Set-TextValue 'B19' 'Uniswap'
Set-TextValue 'C19' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D19' '6.43'
Set-TextValue 'E19' '  +2.97%  '
Set-TextValue 'B20' 'InternetComputer(DFINITY)'
Set-TextValue 'C20' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D20' '12.16'
Set-TextValue 'E20' '  -0.26%  '
Set-TextValue 'D21' '0.0₃0902'
Set-TextValue 'E21' '  +1.50%  '
Set-TextValue 'D22' '68.49'
Set-TextValue 'E22' '  +0.52%  '
Set-TextValue 'D23' '237.85'
Set-TextValue 'E23' '  +1.16%  '
Set-TextValue 'E24' '  +1.04%  '
Set-TextValue 'D25' '2.46'
Set-TextValue 'E25' '  +1.22%  '
Set-TextValue 'E26' '  -0.04%  '
Set-TextValue 'D27' '24.99'
Set-TextValue 'E27' '  +1.96%  '
Set-TextValue 'E28' '  -0.58%  '
Set-TextValue 'D29' '9.43'
Set-TextValue 'E29' '  +3.46%  '
Set-TextValue 'D30' '32.59'
Set-TextValue 'E30' '  +4.02%  '
Set-TextValue 'D31' '0.117'
Set-TextValue 'E31' '  +15.52%  '
Set-TextValue 'E32' '  +7.50%  '
Set-TextValue 'E33' '  +2.01%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  -0.05%  '
Set-TextValue 'D35' '0.0751'
Set-TextValue 'E35' '  +3.72%  '
Set-TextValue 'D36' '133.33'
Set-TextValue 'E36' '  +28.37%  '
Set-TextValue 'E37' '  +3.23%  '
Set-TextValue 'E38' '  +6.57%  '
Set-TextValue 'E39' '  +0.70%  '
Set-TextValue 'D40' '2.28'
Set-TextValue 'E40' '  -1.14%  '
Set-TextValue 'E41' '  +0.05%  '
Set-TextValue 'D42' '21.46'
Set-TextValue 'E42' '  -5.33%  '
Set-TextValue 'D43' '1.946.56'
Set-TextValue 'E44' '  +1.60%  '
Set-TextValue 'E45' '  +2.96%  '
Set-TextValue 'D46' '2.83'
Set-TextValue 'E46' '  +3.54%  '
Set-TextValue 'D47' '9.28'
Set-TextValue 'E47' '  -2.12%  '
Set-TextValue 'D48' '2.635.21'
Set-TextValue 'E48' '  +1.93%  '
Set-TextValue 'D49' '1.57'
Set-TextValue 'E49' '  +5.16%  '
Set-TextValue 'D50' '52.70'
Set-TextValue 'E50' '  -0.20%  '
Set-TextValue 'D51' '72.25'
Set-TextValue 'E51' '  +0.11%  '
